$d = $word.ActiveDocument

$replacements = @(
    @("47×56=2632", "83×59=4897"),
    @("65×89=5785", "81×76=6156"),
    @("41×13=533", "91×60=5460"),
    @("20×48=960", "31×34=1054"),
    @("57×41=2337", "96×15=1440"),
    @("88×85=7480", "39×80=3120"),
    @("25×13=325", "49×33=1617"),
    @("99×53=5247", "36×79=2844"),
    @("59×79=4661", "37×48=1776"),
    @("95×68=6460", "17×14=238"),
    @("64×13=832", "35×98=3430"),
    @("44×93=4092", "27×13=351"),
    @("16×49=784", "61×57=3477"),
    @("45×77=3465", "83×92=7636"),
    @("56×30=1680", "86×57=4902"),
    @("68×96=6528", "76×82=6232"),
    @("63×63=3969", "42×68=2856"),
    @("14×53=742", "29×27=783"),
    @("71×59=4189", "77×70=5390"),
    @("41×18=738", "62×36=2232"),
    @("29×62=1798", "32×17=544"),
    @("87×83=7221", "82×35=2870"),
    @("90×57=5130", "75×36=2700"),
    @("27×59=1593", "58×47=2726"),
    @("52×14=728", "82×41=3362")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
